$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-28
$values = @(
    @(6, 6),
    @(8, 8),
    @(7, 9),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(4, 8),
    @(7, 8),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 6),
    @(1, 5),
    @(1, 2),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(1, 4),
    @(1, 4),
    @(4, 6)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
